$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

$ws.Cells.Item(2,4).Value = 0.03126722457204629
$ws.Cells.Item(2,5).Value = 0.01650570676031604
$ws.Cells.Item(3,4).Value = 0.02886522092845186
$ws.Cells.Item(3,5).Value = -0.008171999184838019
$ws.Cells.Item(4,4).Value = 0.02928660099621385
$ws.Cells.Item(4,5).Value = -0.00992240173005976
$ws.Cells.Item(5,4).Value = 0.06493978618049852
$ws.Cells.Item(5,5).Value = -0.01247935456403781
$ws.Cells.Item(6,4).Value = 0.0158496945729816
$ws.Cells.Item(6,5).Value = 0.006309382422802878
$ws.Cells.Item(7,4).Value = 0.01577910507814944
$ws.Cells.Item(7,5).Value = 0.003951685058156906
$ws.Cells.Item(8,4).Value = 0.02910169573613961
$ws.Cells.Item(8,5).Value = 0.00804495472186284
$ws.Cells.Item(9,4).Value = 0.03463983376958295
$ws.Cells.Item(9,5).Value = -0.002943507302162418
$ws.Cells.Item(10,4).Value = 0.02959935167470635
$ws.Cells.Item(10,5).Value = -0.01411025875432248
$ws.Cells.Item(11,4).Value = 0.03128095030715254
$ws.Cells.Item(11,5).Value = -0.00491443615620879
$ws.Cells.Item(12,4).Value = 0.0121374714725302
$ws.Cells.Item(12,5).Value = -0.06898222940226162
$ws.Cells.Item(13,4).Value = 0.01391789539774141
$ws.Cells.Item(13,5).Value = 0.01282051282051277
$ws.Cells.Item(14,4).Value = 0.01508732802879424
$ws.Cells.Item(14,5).Value = -0.02247088791848628
$ws.Cells.Item(15,4).Value = 0.009052710548364687
$ws.Cells.Item(15,5).Value = 0.007364408248137222
$ws.Cells.Item(16,4).Value = 0.007868963936415337
$ws.Cells.Item(16,5).Value = 0.004485310607759363
$ws.Cells.Item(17,4).Value = 0.02980955150376213
$ws.Cells.Item(17,5).Value = -0.008195966479418004
$ws.Cells.Item(18,4).Value = 0.02571026267331437
$ws.Cells.Item(18,5).Value = -0.0001677852348993536
$ws.Cells.Item(19,4).Value = 0.03093878733914664
$ws.Cells.Item(19,5).Value = 0.03156193554520392
$ws.Cells.Item(20,4).Value = 0.03121232163162128
$ws.Cells.Item(20,5).Value = -0.01049126774720455
$ws.Cells.Item(21,4).Value = 0.04616258839128307
$ws.Cells.Item(21,5).Value = 0.001057661675692767
$ws.Cells.Item(22,4).Value = 0.03429865121122749
$ws.Cells.Item(22,5).Value = 0.0222387377086668
$ws.Cells.Item(23,4).Value = 0.03262509193791497
$ws.Cells.Item(23,5).Value = -0.002163656579619477
$ws.Cells.Item(24,4).Value = 0.03115996775628741
$ws.Cells.Item(24,5).Value = -0.01653735062581396
$ws.Cells.Item(25,4).Value = 0.01482202917738336
$ws.Cells.Item(25,5).Value = -0.01178711751398975
$ws.Cells.Item(26,4).Value = 0.01462967280396572
$ws.Cells.Item(26,5).Value = 0.02782468837957386
$ws.Cells.Item(27,4).Value = 0.03143683544157357
$ws.Cells.Item(27,5).Value = -0.003636363636363549
$ws.Cells.Item(28,4).Value = 0.03038014992032211
$ws.Cells.Item(28,5).Value = -0.02370010843186854
$ws.Cells.Item(29,4).Value = 0.02915228487410266
$ws.Cells.Item(29,5).Value = -0.005327091488760649
$ws.Cells.Item(30,4).Value = 0.02837403569357807
$ws.Cells.Item(30,5).Value = 0.005017103762827713
$ws.Cells.Item(31,4).Value = 0.03376824959033583
$ws.Cells.Item(31,5).Value = 0.007473216618761436
$ws.Cells.Item(32,4).Value = 0.03105075012122766
$ws.Cells.Item(32,5).Value = -0.0004862461794944428
$ws.Cells.Item(33,4).Value = 0.02940601489163826
$ws.Cells.Item(33,5).Value = 0.0001733703190014424
$ws.Cells.Item(34,4).Value = 0.03226253644917978
$ws.Cells.Item(34,5).Value = 0.002820051538872992
$ws.Cells.Item(35,4).Value = 0.03044524912111177
$ws.Cells.Item(35,5).Value = 0.0002318571759796217
$ws.Cells.Item(36,4).Value = 0.03048681849029071
$ws.Cells.Item(36,5).Value = 0.01312065860560852
$ws.Cells.Item(37,4).Value = 0.03319431778096433
$ws.Cells.Item(37,5).Value = 0.003780539672038064
$ws.Cells.Item(38,4).Value = -0.009361121702934683
$ws.Cells.Item(38,5).Value = -0.001116490509928569

$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-05-05 for illustrative purposes only and are subject to change."

$ws.Protect()
